$wb = $excel.ActiveWorkbook

$wsAppBridge = $wb.Worksheets.Item("APPBRIDGE")
$wsEndnoteDown = $wb.Worksheets.Item("ENDNOTE_DOWN")
$wsEndnoteMaint = $wb.Worksheets.Item("ENDNOTE_MAINTENANCE")

# --- Replace placeholder Jira ticket ids with the real ids that were
# --- created for these test cases.

# APPBRIDGE sheet
$wsAppBridge.Range("A2").Value = "OPQA-3493"
$wsAppBridge.Range("A3").Value = "OPQA-3494"
$wsAppBridge.Range("A4").Value = "OPQA-3495"
$wsAppBridge.Range("A6").Value = "OPQA-3496"
$wsAppBridge.Range("A8").Value = "OPQA-3497"
$wsAppBridge.Range("A10").Value = "OPQA-3498"
$wsAppBridge.Range("A11").Value = "OPQA-3499"
$wsAppBridge.Range("A12").Value = "OPQA-3500"
$wsAppBridge.Range("A13").Value = "OPQA-3501"
$wsAppBridge.Range("A14").Value = "OPQA-3502"

# ENDNOTE_DOWN sheet
$wsEndnoteDown.Range("A2").Value = "OPQA-3503"
$wsEndnoteDown.Range("A4").Value = "OPQA-3504"

# ENDNOTE_MAINTENANCE sheet
$wsEndnoteMaint.Range("A2").Value = "OPQA-3505"
$wsEndnoteMaint.Range("A4").Value = "OPQA-3506"

# --- Update the view state: scroll position, selection, and the active
# --- sheet/tab (now ENDNOTE_MAINTENANCE instead of APPBRIDGE).

$wsAppBridge.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$wsAppBridge.Range("L2:L21").Select()

$wsEndnoteDown.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$wsEndnoteDown.Range("M4").Select()

$wsEndnoteMaint.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$wsEndnoteMaint.Range("L4").Select()
